$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# --- Row 11: type changes from VARCHAR2 to decimal ---
$ws.Range("D11").Value = "decimal"

# --- Row 15: type/length/decimal updated and highlighted in yellow ---
$ws.Range("D15").Value = "decimal"
$ws.Range("E15").Value = 14
$ws.Range("F15").Value = 2
$ws.Range("D15").Interior.Color = 65535
$ws.Range("E15").Interior.Color = 65535
$ws.Range("F15").Interior.Color = 65535

# --- Row 20: length updated and highlighted in yellow ---
$ws.Range("E20").Value = 6
$ws.Range("E20").Interior.Color = 65535

# --- Row 22-24: type changes from VARCHAR2 to NVARCHAR2, highlighted ---
$ws.Range("D22").Value = "NVARCHAR2"
$ws.Range("D22").Interior.Color = 65535

$ws.Range("D23").Value = "NVARCHAR2"
$ws.Range("D23").Interior.Color = 65535

$ws.Range("D24").Value = "NVARCHAR2"
$ws.Range("D24").Interior.Color = 65535

# --- Row 25: type changes from VARCHAR2 to DecimalD, highlighted ---
$ws.Range("D25").Value = "DecimalD"
$ws.Range("D25").Interior.Color = 65535

# --- Row 31: type re-set to VARCHAR2 (same displayed text) ---
$ws.Range("D31").Value = "VARCHAR2"

# --- Row 41: add remark "default 1" in column G, matching the style used ---
# used in the neighboring remark cells (e.g. G40)
$ws.Range("G40").Copy()
$ws.Range("G41").PasteSpecial(-4122)
$ws.Range("G41").Value = "default 1"
$excel.CutCopyMode = 0

# --- Update the current selection/view to D11 ---
$ws.Range("D11").Select()
